$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column to Text format so numeric-looking strings
# (e.g. "598.88", "1.00", "0.0810") keep their exact original text
# representation instead of being auto-converted to floating point numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '65.189.49'
$ws.Range("E2").Value = '  -0.16%  '

$ws.Range("D3").Value = '3.552.63'
$ws.Range("E3").Value = '  +0.56%  '

$ws.Range("E4").Value = '  -0.11%  '

$ws.Range("D5").Value = '598.88'
$ws.Range("E5").Value = '  +0.58%  '

$ws.Range("D6").Value = '134.76'
$ws.Range("E6").Value = '  -3.01%  '

$ws.Range("D7").Value = '3.552.61'
$ws.Range("E7").Value = '  +0.61%  '

$ws.Range("E8").Value = '  -0.02%  '

$ws.Range("E9").Value = '  -0.46%  '

$ws.Range("E10").Value = '  -1.98%  '

$ws.Range("D11").Value = '7.03'
$ws.Range("E11").Value = '  -2.47%  '

$ws.Range("D12").Value = '0.387'
$ws.Range("E12").Value = '  -1.09%  '

$ws.Range("D13").Value = '4.154.84'
$ws.Range("E13").Value = '  +0.49%  '

$ws.Range("D14").Value = '0.0000183'
$ws.Range("E14").Value = '  -2.22%  '

$ws.Range("D15").Value = '26.96'
$ws.Range("E15").Value = '  -0.21%  '

$ws.Range("D16").Value = '3.554.48'
$ws.Range("E16").Value = '  +0.79%  '

$ws.Range("E17").Value = '  +0.06%  '

$ws.Range("D18").Value = '65.311.10'

$ws.Range("E19").Value = '  -2.71%  '

$ws.Range("D20").Value = '14.42'
$ws.Range("E20").Value = '  +1.39%  '

$ws.Range("E21").Value = '  -0.38%  '

$ws.Range("D22").Value = '390.49'
$ws.Range("E22").Value = '  -1.21%  '

$ws.Range("D23").Value = '0.579'
$ws.Range("E23").Value = '  +1.53%  '

$ws.Range("D24").Value = '3.695.54'
$ws.Range("E24").Value = '  +0.33%  '

$ws.Range("D25").Value = '74.09'
$ws.Range("E25").Value = '  -0.73%  '

$ws.Range("E26").Value = '  +0.05%  '

$ws.Range("D27").Value = '0.0000113'
$ws.Range("E27").Value = '  -0.66%  '

$ws.Range("D28").Value = '7.79'
$ws.Range("E28").Value = '  +0.41%  '

$ws.Range("B29").Value = 'Fetch.AI'
$ws.Range("C29").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D29").Value = '1.57'
$ws.Range("E29").Value = '  +27.30%  '

$ws.Range("B30").Value = 'Binance-PegBSC-USD'
$ws.Range("C30").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D30").Value = '1.00'
$ws.Range("E30").Value = '  +0.15%  '

$ws.Range("D31").Value = '8.49'
$ws.Range("E31").Value = '  +2.80%  '

$ws.Range("D32").Value = '2.29'
$ws.Range("E32").Value = '  +1.16%  '

$ws.Range("D33").Value = '3.555.12'
$ws.Range("E33").Value = '  +0.10%  '

$ws.Range("D34").Value = '23.97'
$ws.Range("E34").Value = '  +0.83%  '

$ws.Range("E35").Value = '  +0.01%  '

$ws.Range("D36").Value = '0.147'
$ws.Range("E36").Value = '  +0.99%  '

$ws.Range("B37").Value = 'Monero'
$ws.Range("C37").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D37").Value = '169.66'
$ws.Range("E37").Value = '  -0.32%  '

$ws.Range("B38").Value = 'Aptos'
$ws.Range("C38").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D38").Value = '6.93'
$ws.Range("E38").Value = '  -0.81%  '

$ws.Range("E39").Value = '  +0.69%  '

$ws.Range("D40").Value = '5.03'
$ws.Range("E40").Value = '  +2.34%  '

$ws.Range("D41").Value = '0.0810'
$ws.Range("E41").Value = '  +1.38%  '

$ws.Range("D42").Value = '0.825'
$ws.Range("E42").Value = '  +0.61%  '

$ws.Range("D43").Value = '26.28'
$ws.Range("E43").Value = '  -2.12%  '

$ws.Range("D44").Value = '43.02'
$ws.Range("E44").Value = '  +0.62%  '

$ws.Range("D45").Value = '1.25'
$ws.Range("E45").Value = '  +4.45%  '

$ws.Range("D46").Value = '1.00'
$ws.Range("E46").Value = '  -0.19%  '

$ws.Range("E47").Value = '  +0.43%  '

$ws.Range("D48").Value = '1.66'
$ws.Range("E48").Value = '  -0.55%  '

$ws.Range("D49").Value = '2.459.28'
$ws.Range("E49").Value = '  +5.06%  '

$ws.Range("E50").Value = '  +1.42%  '

$ws.Range("D51").Value = '0.0264'
$ws.Range("E51").Value = '  +1.26%  '
